$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 15.24968266666667
$ws.Cells.Item(2, 8).Value = 45.749048
$ws.Cells.Item(2, 9).Value = 0.1107894317754914
$ws.Cells.Item(2, 10).Value = 0.1112362699856105
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 5.450693999999999
$ws.Cells.Item(2, 14).Value = 16.352082
$ws.Cells.Item(2, 15).Value = 0.6387568210835569
$ws.Cells.Item(2, 16).Value = 0.6387568210835569
$ws.Cells.Item(2, 17).Value = 83.12135381310399
$ws.Cells.Item(2, 18).Value = 748.092184317936
$ws.Cells.Item(2, 19).Value = 0.07076750525056649
$ws.Cells.Item(2, 20).Value = 0.0710529262052008

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 15.24968266666667
$ws.Cells.Item(3, 8).Value = 45.749048
$ws.Cells.Item(3, 9).Value = 0.1107894317754914
$ws.Cells.Item(3, 10).Value = 0.1112362699856105
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 3.082591
$ws.Cells.Item(3, 14).Value = 9.247773
$ws.Cells.Item(3, 15).Value = 0.3612431789164431
$ws.Cells.Item(3, 16).Value = 0.3612431789164431
$ws.Cells.Item(3, 17).Value = 47.00853454112267
$ws.Cells.Item(3, 18).Value = 423.0768108701041
$ws.Cells.Item(3, 19).Value = 0.04002192652492491
$ws.Cells.Item(3, 20).Value = 0.04018334378040965

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 24.18733
$ws.Cells.Item(4, 8).Value = 72.56198999999999
$ws.Cells.Item(4, 9).Value = 0.1757217251952191
$ws.Cells.Item(4, 10).Value = 0.1764304496638524
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 5.450693999999999
$ws.Cells.Item(4, 14).Value = 16.352082
$ws.Cells.Item(4, 15).Value = 0.6387568210835569
$ws.Cells.Item(4, 16).Value = 0.6387568210835569
$ws.Cells.Item(4, 17).Value = 131.83773450702
$ws.Cells.Item(4, 18).Value = 1186.53961056318
$ws.Cells.Item(4, 19).Value = 0.1122434505810165
$ws.Cells.Item(4, 20).Value = 0.1126961531696248

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 24.18733
$ws.Cells.Item(5, 8).Value = 72.56198999999999
$ws.Cells.Item(5, 9).Value = 0.1757217251952191
$ws.Cells.Item(5, 10).Value = 0.1764304496638524
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 3.082591
$ws.Cells.Item(5, 14).Value = 9.247773
$ws.Cells.Item(5, 15).Value = 0.3612431789164431
$ws.Cells.Item(5, 16).Value = 0.3612431789164431
$ws.Cells.Item(5, 17).Value = 74.55964577203001
$ws.Cells.Item(5, 18).Value = 671.0368119482699
$ws.Cells.Item(5, 19).Value = 0.06347827461420258
$ws.Cells.Item(5, 20).Value = 0.06373429649422753

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 59.44474733333333
$ws.Cells.Item(6, 8).Value = 178.334242
$ws.Cells.Item(6, 9).Value = 0.4318679885380996
$ws.Cells.Item(6, 10).Value = 0.4336098073732854
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 5.450693999999999
$ws.Cells.Item(6, 14).Value = 16.352082
$ws.Cells.Item(6, 15).Value = 0.6387568210835569
$ws.Cells.Item(6, 16).Value = 0.6387568210835569
$ws.Cells.Item(6, 17).Value = 324.015127621316
$ws.Cells.Item(6, 18).Value = 2916.136148591844
$ws.Cells.Item(6, 19).Value = 0.2758586234863465
$ws.Cells.Item(6, 20).Value = 0.2769712221484132

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 59.44474733333333
$ws.Cells.Item(7, 8).Value = 178.334242
$ws.Cells.Item(7, 9).Value = 0.4318679885380996
$ws.Cells.Item(7, 10).Value = 0.4336098073732854
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 3.082591
$ws.Cells.Item(7, 14).Value = 9.247773
$ws.Cells.Item(7, 15).Value = 0.3612431789164431
$ws.Cells.Item(7, 16).Value = 0.3612431789164431
$ws.Cells.Item(7, 17).Value = 183.2438431270074
$ws.Cells.Item(7, 18).Value = 1649.194588143066
$ws.Cells.Item(7, 19).Value = 0.1560093650517531
$ws.Cells.Item(7, 20).Value = 0.1566385852248722

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 37.10510833333333
$ws.Cells.Item(8, 8).Value = 111.315325
$ws.Cells.Item(8, 9).Value = 0.2695697974885543
$ws.Cells.Item(8, 10).Value = 0.2706570319285326
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 5.450693999999999
$ws.Cells.Item(8, 14).Value = 16.352082
$ws.Cells.Item(8, 15).Value = 0.6387568210835569
$ws.Cells.Item(8, 16).Value = 0.6387568210835569
$ws.Cells.Item(8, 17).Value = 202.2485913618499
$ws.Cells.Item(8, 18).Value = 1820.23732225665
$ws.Cells.Item(8, 19).Value = 0.1721895469039271
$ws.Cells.Item(8, 20).Value = 0.1728840253185802

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 37.10510833333333
$ws.Cells.Item(9, 8).Value = 111.315325
$ws.Cells.Item(9, 9).Value = 0.2695697974885543
$ws.Cells.Item(9, 10).Value = 0.2706570319285326
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 3.082591
$ws.Cells.Item(9, 14).Value = 9.247773
$ws.Cells.Item(9, 15).Value = 0.3612431789164431
$ws.Cells.Item(9, 16).Value = 0.3612431789164431
$ws.Cells.Item(9, 17).Value = 114.3798730023583
$ws.Cells.Item(9, 18).Value = 1029.418857021225
$ws.Cells.Item(9, 19).Value = 0.09738025058462715
$ws.Cells.Item(9, 20).Value = 0.09777300660995236

# Row 10
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 7).Value = 1.6587755
$ws.Cells.Item(10, 8).Value = 3.317551
$ws.Cells.Item(10, 9).Value = 0.01205105700263577
$ws.Cells.Item(10, 10).Value = 0.00806644104871935
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 5.450693999999999
$ws.Cells.Item(10, 14).Value = 16.352082
$ws.Cells.Item(10, 15).Value = 0.6387568210835569
$ws.Cells.Item(10, 16).Value = 0.6387568210835569
$ws.Cells.Item(10, 17).Value = 9.041477665196998
$ws.Cells.Item(10, 18).Value = 54.24886599118199
$ws.Cells.Item(10, 19).Value = 0.00769769486170036
$ws.Cells.Item(10, 20).Value = 0.005152494241737885

# Row 11
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 7).Value = 1.6587755
$ws.Cells.Item(11, 8).Value = 3.317551
$ws.Cells.Item(11, 9).Value = 0.01205105700263577
$ws.Cells.Item(11, 10).Value = 0.00806644104871935
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 3.082591
$ws.Cells.Item(11, 14).Value = 9.247773
$ws.Cells.Item(11, 15).Value = 0.3612431789164431
$ws.Cells.Item(11, 16).Value = 0.3612431789164431
$ws.Cells.Item(11, 17).Value = 5.113326427320501
$ws.Cells.Item(11, 18).Value = 30.679958563923
$ws.Cells.Item(11, 19).Value = 0.004353362140935408
$ws.Cells.Item(11, 20).Value = 0.002913946806981465
